$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7750759720802307
$ws.Range("B1").Value = 1.183954000473022
$ws.Range("C1").Value = 5.270091533660889
$ws.Range("D1").Value = 1.463203549385071
$ws.Range("E1").Value = 0.8542773723602295
